# Fruta / hortaliza, semanal
# The underlying data rows (2..18) get reshuffled: each destination row ends up
# containing the values that used to live in a (different) source row.
# Row 4 is unchanged (maps to itself).
#
# Mapping: destination row -> source row (using the *original* values)
#   2 <- 7    3 <- 14   4 <- 4    5 <- 6    6 <- 16   7 <- 12   8 <- 11
#   9 <- 10   10 <- 5   11 <- 17  12 <- 18  13 <- 8   14 <- 9   15 <- 13
#   16 <- 2   17 <- 15  18 <- 3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: snapshot every data row (columns A..R) BEFORE any writes happen,
# so overwriting one row never destroys data still needed for another row.
$rowData = @{}
for ($r = 2; $r -le 18; $r++) {
    $addr = "A" + $r + ":R" + $r
    $rowData[$r] = $ws.Range($addr).Value2
}

# Step 2: destination row -> source row mapping
$mapping = @{
    2  = 7
    3  = 14
    4  = 4
    5  = 6
    6  = 16
    7  = 12
    8  = 11
    9  = 10
    10 = 5
    11 = 17
    12 = 18
    13 = 8
    14 = 9
    15 = 13
    16 = 2
    17 = 15
    18 = 3
}

# Step 3: write the captured rows back into their new positions.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $addr = "A" + $destRow + ":R" + $destRow
    $ws.Range($addr).Value2 = $rowData[$srcRow]
}
